# Bond dates update: "today" advanced by one day (2023-09-30 -> 2023-10-01).
# Column G ("Dni od poprzedniej wyplaty" = days since previous payment) = TODAY() - F
# Column I ("Dni do nastepnej wyplaty" = days until next payment)      = H - TODAY()
# Since TODAY() moved forward by one day, every existing G value increases by 1
# and every existing I value decreases by 1. Columns F and H (the actual dates)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 262 }

for ($r = 2; $r -le $lastRow; $r++) {
    $gCell = $ws.Cells.Item($r, 7)   # column G
    $gVal = $gCell.Value2
    if ($gVal -ne $null) {
        $gCell.Value2 = $gVal + 1
    }

    $iCell = $ws.Cells.Item($r, 9)   # column I
    $iVal = $iCell.Value2
    if ($iVal -ne $null) {
        $iCell.Value2 = $iVal - 1
    }
}
